$d = $word.ActiveDocument

$pairs = @(
  ,@("2023-04-27 Thursday", "2023-04-28 Friday")
  ,@("21+27=48", "75+20=95")
  ,@("6+20=26", "96-92=4")
  ,@("85+4=89", "79-47=32")
  ,@("77-41=36", "70-1=69")
  ,@("35+18=53", "31-11=20")
  ,@("16+66=82", "38-13=25")
  ,@("2+89=91", "88-64=24")
  ,@("51+14=65", "83-62=21")
  ,@("78+6=84", "92-73=19")
  ,@("52+0=52", "47-9=38")
  ,@("16+79=95", "14-5=9")
  ,@("92-16=76", "18+36=54")
  ,@("95-17=78", "44+12=56")
  ,@("45-41=4", "71-54=17")
  ,@("41+2=43", "59-37=22")
  ,@("10+45=55", "64-51=13")
  ,@("17-4=13", "33+20=53")
  ,@("32-27=5", "99-54=45")
  ,@("0+94=94", "31+32=63")
  ,@("28-12=16", "9+37=46")
  ,@("82-47=35", "3+39=42")
  ,@("4+79=83", "22+4=26")
  ,@("43+56=99", "56+5=61")
  ,@("57-20=37", "71-27=44")
  ,@("12+35=47", "62-7=55")
  ,@("10+18=28", "54+44=98")
  ,@("14+52=66", "40-29=11")
  ,@("58+38=96", "47+37=84")
  ,@("87-23=64", "34+40=74")
  ,@("79-23=56", "62-16=46")
  ,@("19+32=51", "5+81=86")
  ,@("8+55=63", "43-18=25")
  ,@("83-79=4", "71+24=95")
  ,@("54+45=99", "87-73=14")
  ,@("4+22=26", "11-10=1")
  ,@("71-31=40", "16+65=81")
  ,@("23+0=23", "65+2=67")
  ,@("5+44=49", "5+6=11")
  ,@("62+7=69", "92-52=40")
  ,@("26+49=75", "41-2=39")
  ,@("31+67=98", "31+0=31")
  ,@("80+16=96", "52+39=91")
  ,@("37+57=94", "38-34=4")
  ,@("58-43=15", "82-32=50")
  ,@("33+42=75", "9+6=15")
  ,@("35+60=95", "98-4=94")
  ,@("24-3=21", "74-39=35")
  ,@("63-5=58", "24-20=4")
  ,@("84-71=13", "18+75=93")
  ,@("57-1=56", "15+71=86")
  ,@("87+5=92", "79-42=37")
  ,@("91-41=50", "24+8=32")
  ,@("7+85=92", "12+38=50")
  ,@("58+21=79", "37+20=57")
  ,@("96-88=8", "20+7=27")
  ,@("81-52=29", "1+1=2")
  ,@("31-20=11", "74-14=60")
  ,@("26-24=2", "80-66=14")
  ,@("48-36=12", "2+68=70")
  ,@("11-7=4", "20+16=36")
  ,@("6+57=63", "37-19=18")
  ,@("8+39=47", "50+27=77")
  ,@("59+15=74", "50+28=78")
  ,@("10+77=87", "20-5=15")
  ,@("60-28=32", "59-49=10")
  ,@("86-66=20", "91-46=45")
  ,@("94-33=61", "8+30=38")
  ,@("12+79=91", "21+15=36")
  ,@("51+18=69", "68+15=83")
  ,@("94-86=8", "94-38=56")
  ,@("21-17=4", "94-62=32")
  ,@("28-11=17", "14+64=78")
  ,@("37+0=37", "7+39=46")
  ,@("1+8=9", "6+43=49")
  ,@("61+25=86", "54-28=26")
  ,@("84-61=23", "62-60=2")
  ,@("39+5=44", "72+2=74")
  ,@("31-18=13", "95-30=65")
  ,@("36+14=50", "9+2=11")
  ,@("8+63=71", "50+16=66")
  ,@("42-14=28", "4+17=21")
  ,@("73-27=46", "35+32=67")
  ,@("90-0=90", "50-17=33")
  ,@("28+56=84", "30+34=64")
  ,@("57-49=8", "39+36=75")
  ,@("40+1=41", "33+20=53")
  ,@("70+15=85", "40-25=15")
  ,@("31-8=23", "72+23=95")
  ,@("76-46=30", "76-32=44")
  ,@("48-1=47", "41+14=55")
  ,@("40+15=55", "49-42=7")
  ,@("14+75=89", "85-47=38")
  ,@("60-9=51", "5+7=12")
  ,@("8+69=77", "6+83=89")
  ,@("53+46=99", "23+8=31")
  ,@("48-29=19", "36-17=19")
  ,@("32+60=92", "19-17=2")
  ,@("40+58=98", "82-62=20")
  ,@("27+15=42", "95-19=76")
  ,@("46+2=48", "54+3=57")
)

foreach ($pair in $pairs) {
  $old = $pair[0]
  $new = $pair[1]
  $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
